$wb = $excel.ActiveWorkbook

# Suppress the "permanently delete this sheet" confirmation dialog
$excel.DisplayAlerts = $false

# Drop the empty "2019robot" sheet and keep "2018robot" (which holds the data)
$wb.Worksheets.Item("2019robot").Delete()

# Rename the remaining sheet to "Feuil1"
$wb.Worksheets.Item("2018robot").Name = "Feuil1"

$excel.DisplayAlerts = $true
